$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '21 de agosto até 27'
$ws.Range("L2").Value = 'Em cursos na verdade não foi nenhum curso, foi debian, obsidian, itp, C e git'
$ws.Range("B3").Value = '28 até 3 de setembro'
$ws.Range("L3").Value = 'Cursos foi a Proz, e foquei mais em Cálculo pois teoricamente teriamos prova de cálculo na próxima semana'
$ws.Range("B4").Value = '4 até 10 de setembro'
$ws.Range("L4").Value = 'Estudei muito fmc pois tinha lista dela para entregar segunda, e cálculo pois seria prova de Samyr na próxima semana'
$ws.Range("B5").Value = '10 até 17 de setembro'
$ws.Range("L5").Value = 'Tive duas provas na semana (FMC e Cálculo), portanto após elas planejava descansar, teve aniversário de pai fui para Genipabu no FDS e ainda teve as palestras da IEEE na quinta'
$ws.Range("B6").Value = '17 até 24 de setembro'
$ws.Range("L6").Value = 'As questões de itp tavam muito difícil, pré prova do detran, itp, fmc e cansaço'
$ws.Range("B7").Value = '24 até 1 de outubro'
$ws.Range("B8").Value = '1 até 8 de outubro'
$ws.Range("L8").Value = 'Tive duas provas (FMC e Cálculo) e tive as questões mais difíceis da vida em itp(criptografia) tia no domingo'
$ws.Range("B9").Value = '8 até 15 de outubro'
$ws.Range("L9").Value = 'Fiquei doente com distúrbio do sono(sexta) Espírito quebrado cansaço eclipse sábado'
$ws.Range("B10").Value = '15 até 22 de outubro'
$ws.Range("L10").Value = 'Teve todo o negócio da proz de fazer o css da menina lá e o html tbm, e reunião da petcc na quarta além de eu tá dirigindo sábado e domingo'
$ws.Range("B11").Value = '22 até 29 de outubro'
$ws.Range("L11").Value = 'Tive muitos problemas com as listas de itp acumuladas(struct e ponteiros) além de ter prova de FMC na próxima semana + exaustão mental, ademais adicionei a label leituras antes estava contabilizando ela em currículo'
$ws.Range("B12").Value = '29 até 5 de novembro'
$ws.Range("B13").Value = '5 até 12 de novembro'
$ws.Range("L13").Value = 'Pensava que ia ter prova de itp na semana mas o professor adiou, ainda teve prova de cálculo na sexta e sábado e domingo não consegui estudar pq teria prova prática do detran na segundo, assim a parte de leitura é só PLE'
$ws.Range("B14").Value = '12 até 19 de novembro'
$ws.Range("L14").Value = 'Tive a prova do detran na segunda o que me deixou muito ansioso, fui pra igreja também teve uns aniversários e saiu o resultado do detran também o que me fez desopilar'
$ws.Range("B15").Value = '19 até 26 de novembro'
$ws.Range("L15").Value = 'Fui para vó damiana na segunda feira pois pensava que ia ter prova de itp na segunda, mas na realidade foi online, foi uma semana bem complicada mentalmente, eu estava me sentindo muito mal e teve o show do titãs sábado'
$ws.Range("B16").Value = '26 até 3 de dezembro'
$ws.Range("L16").Value = 'Projeto de ITP além de muita coisa pra fazer que me sobrecarregou, meu aniversário foi domingo'
$ws.Range("B17").Value = '3 até 10 de dezembro'
$ws.Range("L17").Value = 'Teve os minicursos da pet na quinta e sexta que comeram muito tempo e eu considerei eles como curriculo, teve o projeto de itp e eu ia ter a prova de fmc na segunda da próxima semana'
